$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation row was inserted as row 114 ("Sin especificar" pea
# variety, Región Metropolitana) pushing all subsequent data rows (the
# former rows 114-164) down by one, down to the new last row 165.
$ws.Rows.Item(114).Insert()

$ws.Cells.Item(114,1).Value() = 6
$ws.Cells.Item(114,2).Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(114,3).Value() = "Metropolitana"
$ws.Cells.Item(114,4).Value() = 44510
$ws.Cells.Item(114,5).Value() = 13
$ws.Cells.Item(114,6).Value() = 100112022
$ws.Cells.Item(114,7).Value() = "Arveja Verde"
$ws.Cells.Item(114,8).Value() = "Sin especificar"
$ws.Cells.Item(114,9).Value() = "Primera"
$ws.Cells.Item(114,10).Value() = 290
$ws.Cells.Item(114,11).Value() = 12000
$ws.Cells.Item(114,12).Value() = 14000
$ws.Cells.Item(114,13).Value() = 13172
$ws.Cells.Item(114,14).Value() = "`$/saco 25 kilos"
$ws.Cells.Item(114,15).Value() = "Región Metropolitana"
$ws.Cells.Item(114,16).Value() = 527
$ws.Cells.Item(114,17).Value() = 25
$ws.Cells.Item(114,18).Value() = "Hortaliza"

# Make sure the date cell uses the same date-formatted number format as
# the rest of column D (the row Insert operation normally already
# copies this from the row above, but we set it explicitly to be safe).
$ws.Cells.Item(114,4).NumberFormat = $ws.Cells.Item(113,4).NumberFormat()
